$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C9 previously held "NA"; the script run now found nothing there, so the
# "NA" marker is cleared, leaving an empty (but still present) text cell -
# same shape as the other blank cells in column C (C2:C8).
$ws.Range("C9").Value = "'"
$ws.Range("C9").ClearFormats()

# Append the new result row produced by the latest script run.
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "2025-02-24"
$ws.Range("A10").ClearFormats()

$ws.Range("B10").Value = "Rien ne nous concerne aujourd'hui !"
$ws.Range("C10").Value = "NA"
$ws.Range("D10").Value = 234
